$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 279.83334
$ws.Range("I5").Value = 59.666668
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 59.666668
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 55.333332
$ws.Range("N5").Value = -730

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7181.4707
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 7181.4707
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 21544.4121
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -21880.4121

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 20834288
$ws.Range("I19").Value = 66667130
$ws.Range("J19").Value = 1180.3636
$ws.Range("K19").Value = 66667130
$ws.Range("L19").Value = 1180.3636
$ws.Range("M19").Value = -66666955
$ws.Range("N19").Value = -1530.3636

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1955
$ws.Range("J97").Value = 1955
$ws.Range("L97").Value = 5865
$ws.Range("N97").Value = -6857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2824.6875
$ws.Range("I113").Value = 1784
$ws.Range("K113").Value = 1784
$ws.Range("M113").Value = 1470

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 982.3
$ws.Range("J129").Value = 1010.7292
$ws.Range("L129").Value = 3032.1876
$ws.Range("N129").Value = -13032.1876

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5738.213
$ws.Range("I32").Value = 4448.0186
$ws.Range("K32").Value = 4448.0186
$ws.Range("M32").Value = -4161.0186

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2019.125
$ws.Range("I45").Value = 1869.0526
$ws.Range("J45").Value = 2238.4614
$ws.Range("K45").Value = 1869.0526
$ws.Range("L45").Value = 2238.4614
$ws.Range("M45").Value = -1492.0526
$ws.Range("N45").Value = -2992.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1123.069
$ws.Range("I97").Value = 849.9048
$ws.Range("J97").Value = 1840.125
$ws.Range("K97").Value = 849.9048
$ws.Range("L97").Value = 1840.125
$ws.Range("M97").Value = -353.9048
$ws.Range("N97").Value = -2832.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2218.7778
$ws.Range("I102").Value = 2105.3845
$ws.Range("J102").Value = 2324.0715
$ws.Range("K102").Value = 2105.3845
$ws.Range("L102").Value = 2324.0715
$ws.Range("M102").Value = -483.3845000000001
$ws.Range("N102").Value = -5568.0715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4546.2617
$ws.Range("I132").Value = 1542.8846
$ws.Range("J132").Value = 9426.75
$ws.Range("K132").Value = 4628.6538
$ws.Range("L132").Value = 28280.25
$ws.Range("M132").Value = -2098.6538
$ws.Range("N132").Value = -33340.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1563.4333
$ws.Range("I94").Value = 1439.6111
$ws.Range("J94").Value = 1749.1666
$ws.Range("K94").Value = 1439.6111
$ws.Range("L94").Value = 1749.1666
$ws.Range("M94").Value = -988.6111000000001
$ws.Range("N94").Value = -2651.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2306.4119
$ws.Range("I99").Value = 2141
$ws.Range("J99").Value = 2542.7144
$ws.Range("K99").Value = 2141
$ws.Range("L99").Value = 2542.7144
$ws.Range("M99").Value = -643
$ws.Range("N99").Value = -5538.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 781
$ws.Range("I16").Value = 559.4545000000001
$ws.Range("J16").Value = 1999.5
$ws.Range("K16").Value = 559.4545000000001
$ws.Range("L16").Value = 1999.5
$ws.Range("M16").Value = -272.4545000000001
$ws.Range("N16").Value = -2573.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5708.3335
$ws.Range("I31").Value = 5739.5356
$ws.Range("J31").Value = 5599.125
$ws.Range("K31").Value = 5739.5356
$ws.Range("L31").Value = 5599.125
$ws.Range("M31").Value = -5444.5356
$ws.Range("N31").Value = -6189.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5708.3335
$ws.Range("I34").Value = 5739.5356
$ws.Range("J34").Value = 5599.125
$ws.Range("K34").Value = 5739.5356
$ws.Range("L34").Value = 5599.125
$ws.Range("M34").Value = -5537.5356
$ws.Range("N34").Value = -6003.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 781
$ws.Range("I113").Value = 559.4545000000001
$ws.Range("J113").Value = 1999.5
$ws.Range("K113").Value = 559.4545000000001
$ws.Range("L113").Value = 1999.5
$ws.Range("M113").Value = 1610.5455
$ws.Range("N113").Value = -6339.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 864.0625
$ws.Range("I26").Value = 100.5
$ws.Range("J26").Value = 1322.2
$ws.Range("K26").Value = 301.5
$ws.Range("L26").Value = 3966.6
$ws.Range("M26").Value = -13.5
$ws.Range("N26").Value = -4542.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1640234.1
$ws.Range("I139").Value = 3355491.2
$ws.Range("J139").Value = 2943.0454
$ws.Range("K139").Value = 10066473.6
$ws.Range("L139").Value = 8829.136200000001
$ws.Range("M139").Value = -10061333.6
$ws.Range("N139").Value = -19109.1362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 77008460
$ws.Range("I14").Value = 77008460
$ws.Range("K14").Value = 77008460
$ws.Range("M14").Value = -77008292

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1635.5172
$ws.Range("I132").Value = 1134.3334
$ws.Range("J132").Value = 2455.6365
$ws.Range("K132").Value = 3403.0002
$ws.Range("L132").Value = 7366.9095
$ws.Range("M132").Value = -873.0001999999999
$ws.Range("N132").Value = -12426.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 764714.4399999999
$ws.Range("I2").Value = 4000.5
$ws.Range("J2").Value = 1069000
$ws.Range("K2").Value = 4000.5
$ws.Range("L2").Value = 1069000
$ws.Range("M2").Value = -3888.5
$ws.Range("N2").Value = -1069224

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 419.23077
$ws.Range("I22").Value = 325.16666
$ws.Range("J22").Value = 499.85715
$ws.Range("K22").Value = 325.16666
$ws.Range("L22").Value = 499.85715
$ws.Range("M22").Value = -30.16665999999998
$ws.Range("N22").Value = -1089.85715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 419.23077
$ws.Range("I27").Value = 325.16666
$ws.Range("J27").Value = 499.85715
$ws.Range("K27").Value = 325.16666
$ws.Range("L27").Value = 499.85715
$ws.Range("M27").Value = -218.16666
$ws.Range("N27").Value = -713.85715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 963.5
$ws.Range("I93").Value = 963.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 963.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 284.5
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3020.348
$ws.Range("I132").Value = 2234.2068
$ws.Range("J132").Value = 4361.4116
$ws.Range("K132").Value = 6702.6204
$ws.Range("L132").Value = 13084.2348
$ws.Range("M132").Value = -4172.6204
$ws.Range("N132").Value = -18144.2348

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1430.5
$ws.Range("I96").Value = 1137.5
$ws.Range("J96").Value = 1723.5
$ws.Range("K96").Value = 1137.5
$ws.Range("L96").Value = 1723.5
$ws.Range("M96").Value = 235.5
$ws.Range("N96").Value = -4469.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1711.4
$ws.Range("I100").Value = 786.8570999999999
$ws.Range("K100").Value = 1573.7142
$ws.Range("M100").Value = -1032.7142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 6501.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 6501.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 19504.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -23344.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1894.4872
$ws.Range("I122").Value = 1684.5588
$ws.Range("K122").Value = 5053.6764
$ws.Range("M122").Value = -2603.6764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1931.9667
$ws.Range("I132").Value = 1014.44446
$ws.Range("J132").Value = 3308.25
$ws.Range("K132").Value = 3043.33338
$ws.Range("L132").Value = 9924.75
$ws.Range("M132").Value = -513.33338
$ws.Range("N132").Value = -14984.75
